$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 = "Save", styled like the other header cells (copy style from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column H data values (Save) for rows 2-11
$values = @(0, 0, 0, 0, 0, 0, 0, 1, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
